$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new "DRINKS" section rows (76-86) ---
$ws.Range('B76').Value = 'DRINKS'
$ws.Range('C76').Value = 'N/A'
$ws.Range('D76').Value = 'NAPOJE'
$ws.Range('E76').Value = 'N/A'
$ws.Range('F76').Value = 'N/A'
$ws.Range('G76').Value = 'N/A'
$ws.Range('H76').Value = 'N/A'
$ws.Range('I76').Value = 'N/A'
$ws.Range('F77').Value = 'PEPSI/PEPSI MAX'
$ws.Range('F78').Value = 'MIRINDA'
$ws.Range('F79').Value = '7UP'
$ws.Range('F80').Value = 'LIPTON'
$ws.Range('F82').Value = 'HOMEMADE LEMONADE'
$ws.Range('F81').Value = 'TOMA JUICE'
$ws.Range('F83').Value = 'BASIL SEEDS JUICE'
$ws.Range('F84').Value = 'WATER'
$ws.Range('F85').Value = 'BUBBLE TEA'
$ws.Range('F86').Value = 'GUAVA JUICE'
$ws.Range('G77').Value = 'PEPSI/PEPSI MAX'
$ws.Range('G78').Value = 'MIRINDA'
$ws.Range('G79').Value = '7UP'
$ws.Range('G80').Value = 'LIPTON'
$ws.Range('G81').Value = 'TOMA JUICE'
$ws.Range('G82').Value = 'DOMOWA LEMONIADA'
$ws.Range('G83').Value = 'SOK Z PESTKAMI BAZYLII'
$ws.Range('G84').Value = 'WODA'
$ws.Range('G85').Value = 'BUBBLE TEA'
$ws.Range('G86').Value = 'SOK Z GUAWY'
$ws.Range('H77').Value = '200 ML/500 ML'
$ws.Range('H78').Value = '200 ML/500 ML'
$ws.Range('H79').Value = '200 ML/500 ML'
$ws.Range('H80').Value = '200 ML/500 ML'
$ws.Range('H81').Value = '200 ML/1 L'
$ws.Range('H82').Value = '400 ML/ 1 L'
$ws.Range('H84').Value = '500 ML/1 L'
$ws.Range('I77').Value = '8 PLN/10 PLN'
$ws.Range('I78').Value = '8 PLN/10 PLN'
$ws.Range('I79').Value = '8 PLN/10 PLN'
$ws.Range('I80').Value = '8 PLN/10 PLN'
$ws.Range('I81').Value = '8 PLN/18 PLN'
$ws.Range('I82').Value = '16 PLN/25 PLN'
$ws.Range('I84').Value = '8 PLN/16 PLN'
$ws.Range('I83').Value = '16 PLN'
$ws.Range('I85').Value = '16 PLN'
$ws.Range('I86').Value = '16 PLN'
$ws.Range('B77').Value = 'DRINKS'
$ws.Range('C77').Value = 'N/A'
$ws.Range('D77').Value = 'NAPOJE'
$ws.Range('E77').Value = 'N/A'
$ws.Range('B78').Value = 'DRINKS'
$ws.Range('C78').Value = 'N/A'
$ws.Range('D78').Value = 'NAPOJE'
$ws.Range('E78').Value = 'N/A'
$ws.Range('B79').Value = 'DRINKS'
$ws.Range('C79').Value = 'N/A'
$ws.Range('D79').Value = 'NAPOJE'
$ws.Range('E79').Value = 'N/A'
$ws.Range('B80').Value = 'DRINKS'
$ws.Range('C80').Value = 'N/A'
$ws.Range('D80').Value = 'NAPOJE'
$ws.Range('E80').Value = 'N/A'
$ws.Range('B81').Value = 'DRINKS'
$ws.Range('C81').Value = 'N/A'
$ws.Range('D81').Value = 'NAPOJE'
$ws.Range('E81').Value = 'N/A'
$ws.Range('B82').Value = 'DRINKS'
$ws.Range('C82').Value = 'N/A'
$ws.Range('D82').Value = 'NAPOJE'
$ws.Range('E82').Value = 'N/A'
$ws.Range('B83').Value = 'DRINKS'
$ws.Range('C83').Value = 'N/A'
$ws.Range('D83').Value = 'NAPOJE'
$ws.Range('E83').Value = 'N/A'
$ws.Range('B84').Value = 'DRINKS'
$ws.Range('C84').Value = 'N/A'
$ws.Range('D84').Value = 'NAPOJE'
$ws.Range('E84').Value = 'N/A'
$ws.Range('B85').Value = 'DRINKS'
$ws.Range('C85').Value = 'N/A'
$ws.Range('D85').Value = 'NAPOJE'
$ws.Range('E85').Value = 'N/A'
$ws.Range('B86').Value = 'DRINKS'
$ws.Range('C86').Value = 'N/A'
$ws.Range('D86').Value = 'NAPOJE'
$ws.Range('E86').Value = 'N/A'

# --- Bold the section/sub-section header cells on the first row of the new block ---
$ws.Range('B76').Font.Bold = $true
$ws.Range('D76').Font.Bold = $true

# --- Restore the active selection left by the editor ---
$ws.Range('D69').Select()
